# Fruta / hortaliza, semanal
# The underlying data rows (2..39, excluding the untouched row 14) got
# reshuffled: each row's Fecha/Volumen/Precio.../Unidad de comercializacion/
# Precio $/Kg / Kg o Unidades values were swapped with those of another row
# in the same column set, while the descriptive columns (Mercado, Region,
# Codreg, Categoria, Variedad, Calidad, Origen, Clasificacion) stayed put
# because they are identical for every row anyway.
#
# Build the row -> source-row mapping observed in the target workbook, then
# snapshot the "before" values and re-write them into their new homes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destinationRow -> sourceRow (values that should end up in destinationRow
# are the ones that used to live in sourceRow)
$rowMap = @{
    2=37; 3=5; 4=35; 5=4; 6=7; 7=38; 8=25; 9=30; 10=22;
    11=32; 12=29; 13=23; 15=20; 16=17; 17=31; 18=2; 19=27; 20=9;
    21=28; 22=13; 23=36; 24=12; 25=18; 26=33; 27=11; 28=6; 29=19;
    30=39; 31=3; 32=34; 33=8; 34=24; 35=16; 36=21; 37=10; 38=15; 39=26
}

# Columns involved in the shuffle, by column index.
$cols = @(4, 10, 11, 12, 13, 14, 16, 17)   # D, J, K, L, M, N, P, Q

# 1) Snapshot every relevant cell value before mutating anything.
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowData
}

# 2) Write each destination row's cells using the snapshot of its source row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcData[$c]
    }
}
